$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
}

Set-TextValue $ws.Range("D2") "304.54"
Set-TextValue $ws.Range("E2") "2.32%"
Set-TextValue $ws.Range("D3") "31.90"
Set-TextValue $ws.Range("E3") "0.62%"
Set-TextValue $ws.Range("E4") "0.37%"
Set-TextValue $ws.Range("D5") "0.07473"
Set-TextValue $ws.Range("E5") "-0.01%"
Set-TextValue $ws.Range("D6") "2.397"
Set-TextValue $ws.Range("E6") "42.15%"
Set-TextValue $ws.Range("D7") "8.009"
Set-TextValue $ws.Range("E7") "3.15%"
$ws.Range("B8").Value2 = "MXToken"
$ws.Range("C8").Value2 = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.9175"
Set-TextValue $ws.Range("E8") "-0.81%"
$ws.Range("B9").Value2 = "WazirX"
$ws.Range("C9").Value2 = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D9") "0.1734"
Set-TextValue $ws.Range("E9") "1.54%"
$ws.Range("B10").Value2 = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value2 = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.07671"
Set-TextValue $ws.Range("E10") "2.21%"
$ws.Range("B11").Value2 = "MandalaExchangeToken"
$ws.Range("C11").Value2 = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.08163"
Set-TextValue $ws.Range("E11") "3.12%"
$ws.Range("B12").Value2 = "BitrueCoin"
$ws.Range("C12").Value2 = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D12") "0.03012"
Set-TextValue $ws.Range("E12") "0.56%"
$ws.Range("B13").Value2 = "BitMartToken"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D13") "0.09935"
Set-TextValue $ws.Range("E13") "0.46%"
$ws.Range("B14").Value2 = "BitForexToken"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D14") "0.001506"
Set-TextValue $ws.Range("E14") "-0.24%"
$ws.Range("B15").Value2 = "TigerCash"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D15") "0.006060"
Set-TextValue $ws.Range("E15") "-2.60%"
$ws.Range("B16").Value2 = "LEO"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D16") "3.507"
Set-TextValue $ws.Range("E16") "1.89%"
$ws.Range("B17").Value2 = "GateToken"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D17") "3.859"
Set-TextValue $ws.Range("E17") "1.35%"
Set-TextValue $ws.Range("E18") "-0.09%"
Set-TextValue $ws.Range("E19") "-0.95%"
Set-TextValue $ws.Range("E20") "-0.76%"
Set-TextValue $ws.Range("D21") "4.651"
Set-TextValue $ws.Range("E21") "1.94%"
Set-TextValue $ws.Range("D22") "0.04605"
Set-TextValue $ws.Range("E22") "-1.04%"
Set-TextValue $ws.Range("D23") "0.1563"
Set-TextValue $ws.Range("E23") "0.70%"
Set-TextValue $ws.Range("D24") "0.001260"
Set-TextValue $ws.Range("E24") "3.02%"
Set-TextValue $ws.Range("D25") "0.004529"
Set-TextValue $ws.Range("E25") "2.67%"
Set-TextValue $ws.Range("D26") "0.0001299"
Set-TextValue $ws.Range("E26") "-7.30%"
Set-TextValue $ws.Range("D27") "0.0002737"
Set-TextValue $ws.Range("E27") "51.12%"
Set-TextValue $ws.Range("D39") "0.01777"
Set-TextValue $ws.Range("E39") "7.71%"
Set-TextValue $ws.Range("D40") "0.04566"
Set-TextValue $ws.Range("E40") "1.22%"
Set-TextValue $ws.Range("D41") "0.007436"
Set-TextValue $ws.Range("E41") "7.14%"
Set-TextValue $ws.Range("D42") "0.1362"
Set-TextValue $ws.Range("E42") "1.48%"
Set-TextValue $ws.Range("D43") "0.002178"
Set-TextValue $ws.Range("E43") "5.65%"
Set-TextValue $ws.Range("D44") "0.01079"
Set-TextValue $ws.Range("E44") "-18.60%"
Set-TextValue $ws.Range("D45") "0.00006435"
Set-TextValue $ws.Range("E45") "5.56%"
Set-TextValue $ws.Range("E46") "-57.48%"
Set-TextValue $ws.Range("D47") "0.009885"
Set-TextValue $ws.Range("E47") "-19.35%"
